$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 27.67568533333333
$ws.Range("H2").Value = 83.027056
$ws.Range("I2").Value = 0.151580065893459
$ws.Range("J2").Value = 0.151580065893459
$ws.Range("M2").Value = 14.89002333333333
$ws.Range("N2").Value = 44.67007
$ws.Range("O2").Value = 0.1194491234330596
$ws.Range("P2").Value = 0.1194491234330597
$ws.Range("Q2").Value = 412.0916003793245
$ws.Range("R2").Value = 3708.82440341392
$ws.Range("S2").Value = 0.01810610600089909
$ws.Range("T2").Value = 0.0181061060008991
$ws.Range("G3").Value = 27.67568533333333
$ws.Range("H3").Value = 83.027056
$ws.Range("I3").Value = 0.151580065893459
$ws.Range("J3").Value = 0.151580065893459
$ws.Range("O3").Value = 0.05148509068166413
$ws.Range("P3").Value = 0.05148509068166414
$ws.Range("Q3").Value = 177.6201683603947
$ws.Range("R3").Value = 1598.581515243552
$ws.Range("S3").Value = 0.007804113438057359
$ws.Range("T3").Value = 0.007804113438057362
$ws.Range("G4").Value = 27.67568533333333
$ws.Range("H4").Value = 83.027056
$ws.Range("I4").Value = 0.151580065893459
$ws.Range("J4").Value = 0.151580065893459
$ws.Range("M4").Value = 66.36284166666667
$ws.Range("N4").Value = 199.088525
$ws.Range("O4").Value = 0.5323687604884161
$ws.Range("P4").Value = 0.5323687604884162
$ws.Range("Q4").Value = 1836.637123792489
$ws.Range("R4").Value = 16529.7341141324
$ws.Range("S4").Value = 0.08069649179445318
$ws.Range("T4").Value = 0.08069649179445321
$ws.Range("G5").Value = 27.67568533333333
$ws.Range("H5").Value = 83.027056
$ws.Range("I5").Value = 0.151580065893459
$ws.Range("J5").Value = 0.151580065893459
$ws.Range("M5").Value = 3.521285666666667
$ws.Range("N5").Value = 10.563857
$ws.Range("O5").Value = 0.02824807435318976
$ws.Range("P5").Value = 0.02824807435318976
$ws.Range("Q5").Value = 97.45399407944356
$ws.Range("R5").Value = 877.085946714992
$ws.Range("S5").Value = 0.004281844971819832
$ws.Range("T5").Value = 0.004281844971819833
$ws.Range("G6").Value = 27.67568533333333
$ws.Range("H6").Value = 83.027056
$ws.Range("I6").Value = 0.151580065893459
$ws.Range("J6").Value = 0.151580065893459
$ws.Range("M6").Value = 33.46371266666667
$ws.Range("N6").Value = 100.391138
$ws.Range("O6").Value = 0.2684489510436703
$ws.Range("P6").Value = 0.2684489510436703
$ws.Range("Q6").Value = 926.1311818477476
$ws.Range("R6").Value = 8335.180636629728
$ws.Range("S6").Value = 0.04069150968822948
$ws.Range("T6").Value = 0.04069150968822949
$ws.Range("I7").Value = 0.2439851776203359
$ws.Range("J7").Value = 0.243985177620336
$ws.Range("M7").Value = 14.89002333333333
$ws.Range("N7").Value = 44.67007
$ws.Range("O7").Value = 0.1194491234330596
$ws.Range("P7").Value = 0.1194491234330597
$ws.Range("Q7").Value = 663.3078150597156
$ws.Range("R7").Value = 5969.770335537441
$ws.Range("S7").Value = 0.02914381559740849
$ws.Range("T7").Value = 0.0291438155974085
$ws.Range("I8").Value = 0.2439851776203359
$ws.Range("J8").Value = 0.243985177620336
$ws.Range("O8").Value = 0.05148509068166413
$ws.Range("P8").Value = 0.05148509068166414
$ws.Range("S8").Value = 0.01256159899476493
$ws.Range("T8").Value = 0.01256159899476493
$ws.Range("I9").Value = 0.2439851776203359
$ws.Range("J9").Value = 0.243985177620336
$ws.Range("M9").Value = 66.36284166666667
$ws.Range("N9").Value = 199.088525
$ws.Range("O9").Value = 0.5323687604884161
$ws.Range("P9").Value = 0.5323687604884162
$ws.Range("Q9").Value = 2956.274179136311
$ws.Range("R9").Value = 26606.4676122268
$ws.Range("S9").Value = 0.1298900865872843
$ws.Range("T9").Value = 0.1298900865872843
$ws.Range("I10").Value = 0.2439851776203359
$ws.Range("J10").Value = 0.243985177620336
$ws.Range("M10").Value = 3.521285666666667
$ws.Range("N10").Value = 10.563857
$ws.Range("O10").Value = 0.02824807435318976
$ws.Range("P10").Value = 0.02824807435318976
$ws.Range("Q10").Value = 156.8631727076605
$ws.Range("R10").Value = 1411.768554368944
$ws.Range("S10").Value = 0.00689211143849546
$ws.Range("T10").Value = 0.006892111438495461
$ws.Range("I11").Value = 0.2439851776203359
$ws.Range("J11").Value = 0.243985177620336
$ws.Range("M11").Value = 33.46371266666667
$ws.Range("N11").Value = 100.391138
$ws.Range("O11").Value = 0.2684489510436703
$ws.Range("P11").Value = 0.2684489510436703
$ws.Range("Q11").Value = 1490.712380753788
$ws.Range("R11").Value = 13416.4114267841
$ws.Range("S11").Value = 0.06549756500238277
$ws.Range("T11").Value = 0.06549756500238277
$ws.Range("G12").Value = 54.059897
$ws.Range("H12").Value = 162.179691
$ws.Range("I12").Value = 0.2960867147735651
$ws.Range("J12").Value = 0.2960867147735651
$ws.Range("M12").Value = 14.89002333333333
$ws.Range("N12").Value = 44.67007
$ws.Range("O12").Value = 0.1194491234330596
$ws.Range("P12").Value = 0.1194491234330597
$ws.Range("Q12").Value = 804.9531277275967
$ws.Range("R12").Value = 7244.57814954837
$ws.Range("S12").Value = 0.0353672985398767
$ws.Range("T12").Value = 0.0353672985398767
$ws.Range("G13").Value = 54.059897
$ws.Range("H13").Value = 162.179691
$ws.Range("I13").Value = 0.2960867147735651
$ws.Range("J13").Value = 0.2960867147735651
$ws.Range("O13").Value = 0.05148509068166413
$ws.Range("P13").Value = 0.05148509068166414
$ws.Range("Q13").Value = 346.951769794858
$ws.Range("R13").Value = 3122.565928153722
$ws.Range("S13").Value = 0.01524405135975302
$ws.Range("T13").Value = 0.01524405135975302
$ws.Range("G14").Value = 54.059897
$ws.Range("H14").Value = 162.179691
$ws.Range("I14").Value = 0.2960867147735651
$ws.Range("J14").Value = 0.2960867147735651
$ws.Range("M14").Value = 66.36284166666667
$ws.Range("N14").Value = 199.088525
$ws.Range("O14").Value = 0.5323687604884161
$ws.Range("P14").Value = 0.5323687604884162
$ws.Range("Q14").Value = 3587.568385127308
$ws.Range("R14").Value = 32288.11546614577
$ws.Range("S14").Value = 0.15762731734109
$ws.Range("T14").Value = 0.1576273173410901
$ws.Range("G15").Value = 54.059897
$ws.Range("H15").Value = 162.179691
$ws.Range("I15").Value = 0.2960867147735651
$ws.Range("J15").Value = 0.2960867147735651
$ws.Range("M15").Value = 3.521285666666667
$ws.Range("N15").Value = 10.563857
$ws.Range("O15").Value = 0.02824807435318976
$ws.Range("P15").Value = 0.02824807435318976
$ws.Range("Q15").Value = 190.3603404475763
$ws.Range("R15").Value = 1713.243064028187
$ws.Range("S15").Value = 0.008363879533915356
$ws.Range("T15").Value = 0.008363879533915356
$ws.Range("G16").Value = 54.059897
$ws.Range("H16").Value = 162.179691
$ws.Range("I16").Value = 0.2960867147735651
$ws.Range("J16").Value = 0.2960867147735651
$ws.Range("M16").Value = 33.46371266666667
$ws.Range("N16").Value = 100.391138
$ws.Range("O16").Value = 0.2684489510436703
$ws.Range("P16").Value = 0.2684489510436703
$ws.Range("Q16").Value = 1809.044859997595
$ws.Range("R16").Value = 16281.40373997836
$ws.Range("S16").Value = 0.07948416799892993
$ws.Range("T16").Value = 0.07948416799892993
$ws.Range("G17").Value = 11.41370466666667
$ws.Range("H17").Value = 34.241114
$ws.Range("I17").Value = 0.0625129995743248
$ws.Range("J17").Value = 0.0625129995743248
$ws.Range("M17").Value = 14.89002333333333
$ws.Range("N17").Value = 44.67007
$ws.Range("O17").Value = 0.1194491234330596
$ws.Range("P17").Value = 0.1194491234330597
$ws.Range("Q17").Value = 169.9503288064422
$ws.Range("R17").Value = 1529.55295925798
$ws.Range("S17").Value = 0.007467123002324327
$ws.Range("T17").Value = 0.007467123002324328
$ws.Range("G18").Value = 11.41370466666667
$ws.Range("H18").Value = 34.241114
$ws.Range("I18").Value = 0.0625129995743248
$ws.Range("J18").Value = 0.0625129995743248
$ws.Range("O18").Value = 0.05148509068166413
$ws.Range("P18").Value = 0.05148509068166414
$ws.Range("Q18").Value = 73.25217497206533
$ws.Range("R18").Value = 659.2695747485881
$ws.Range("S18").Value = 0.003218487451866944
$ws.Range("T18").Value = 0.003218487451866944
$ws.Range("G19").Value = 11.41370466666667
$ws.Range("H19").Value = 34.241114
$ws.Range("I19").Value = 0.0625129995743248
$ws.Range("J19").Value = 0.0625129995743248
$ws.Range("M19").Value = 66.36284166666667
$ws.Range("N19").Value = 199.088525
$ws.Range("O19").Value = 0.5323687604884161
$ws.Range("P19").Value = 0.5323687604884162
$ws.Range("Q19").Value = 757.4458756240946
$ws.Range("R19").Value = 6817.01288061685
$ws.Range("S19").Value = 0.03327996809779617
$ws.Range("T19").Value = 0.03327996809779618
$ws.Range("G20").Value = 11.41370466666667
$ws.Range("H20").Value = 34.241114
$ws.Range("I20").Value = 0.0625129995743248
$ws.Range("J20").Value = 0.0625129995743248
$ws.Range("M20").Value = 3.521285666666667
$ws.Range("N20").Value = 10.563857
$ws.Range("O20").Value = 0.02824807435318976
$ws.Range("P20").Value = 0.02824807435318976
$ws.Range("Q20").Value = 40.19091464629979
$ws.Range("R20").Value = 361.7182318166981
$ws.Range("S20").Value = 0.001765871860016447
$ws.Range("T20").Value = 0.001765871860016447
$ws.Range("G21").Value = 11.41370466666667
$ws.Range("H21").Value = 34.241114
$ws.Range("I21").Value = 0.0625129995743248
$ws.Range("J21").Value = 0.0625129995743248
$ws.Range("M21").Value = 33.46371266666667
$ws.Range("N21").Value = 100.391138
$ws.Range("O21").Value = 0.2684489510436703
$ws.Range("P21").Value = 0.2684489510436703
$ws.Range("Q21").Value = 381.9449334275258
$ws.Range("R21").Value = 3437.504400847732
$ws.Range("S21").Value = 0.0167815491623209
$ws.Range("T21").Value = 0.0167815491623209
$ws.Range("G22").Value = 44.88488133333333
$ws.Range("H22").Value = 134.654644
$ws.Range("I22").Value = 0.2458350421383152
$ws.Range("J22").Value = 0.2458350421383153
$ws.Range("M22").Value = 14.89002333333333
$ws.Range("N22").Value = 44.67007
$ws.Range("O22").Value = 0.1194491234330596
$ws.Range("P22").Value = 0.1194491234330597
$ws.Range("Q22").Value = 668.3369303672312
$ws.Range("R22").Value = 6015.03237330508
$ws.Range("S22").Value = 0.02936478029255104
$ws.Range("T22").Value = 0.02936478029255104
$ws.Range("G23").Value = 44.88488133333333
$ws.Range("H23").Value = 134.654644
$ws.Range("I23").Value = 0.2458350421383152
$ws.Range("J23").Value = 0.2458350421383153
$ws.Range("O23").Value = 0.05148509068166413
$ws.Range("P23").Value = 0.05148509068166414
$ws.Range("Q23").Value = 288.0673082975387
$ws.Range("R23").Value = 2592.605774677848
$ws.Range("S23").Value = 0.01265683943722188
$ws.Range("T23").Value = 0.01265683943722189
$ws.Range("G24").Value = 44.88488133333333
$ws.Range("H24").Value = 134.654644
$ws.Range("I24").Value = 0.2458350421383152
$ws.Range("J24").Value = 0.2458350421383153
$ws.Range("M24").Value = 66.36284166666667
$ws.Range("N24").Value = 199.088525
$ws.Range("O24").Value = 0.5323687604884161
$ws.Range("P24").Value = 0.5323687604884162
$ws.Range("Q24").Value = 2978.688273151122
$ws.Range("R24").Value = 26808.1944583601
$ws.Range("S24").Value = 0.1308748966677924
$ws.Range("T24").Value = 0.1308748966677925
$ws.Range("G25").Value = 44.88488133333333
$ws.Range("H25").Value = 134.654644
$ws.Range("I25").Value = 0.2458350421383152
$ws.Range("J25").Value = 0.2458350421383153
$ws.Range("M25").Value = 3.521285666666667
$ws.Range("N25").Value = 10.563857
$ws.Range("O25").Value = 0.02824807435318976
$ws.Range("P25").Value = 0.02824807435318976
$ws.Range("Q25").Value = 158.0524892891009
$ws.Range("R25").Value = 1422.472403601908
$ws.Range("S25").Value = 0.006944366548942666
$ws.Range("T25").Value = 0.006944366548942667
$ws.Range("G26").Value = 44.88488133333333
$ws.Range("H26").Value = 134.654644
$ws.Range("I26").Value = 0.2458350421383152
$ws.Range("J26").Value = 0.2458350421383153
$ws.Range("M26").Value = 33.46371266666667
$ws.Range("N26").Value = 100.391138
$ws.Range("O26").Value = 0.2684489510436703
$ws.Range("P26").Value = 0.2684489510436703
$ws.Range("Q26").Value = 1502.014772016097
$ws.Range("R26").Value = 13518.13294814487
$ws.Range("S26").Value = 0.06599415919180721
$ws.Range("T26").Value = 0.06599415919180721
